# SEIRS.xlsx bugfix:
#  - pars sheet: transmission rate constant default 0.8 -> 0.2
#                incubation rate constant default 0.5 -> 0.4
#  - eqns sheet: stoichiometry columns (S,E,I,R) for the transmission,
#    infection and recovery reactions were mistakenly stored as the
#    string names of the rate-constant parameters / their negation
#    (e.g. "-t", "t", "-i", "i", "-r", "r") instead of the numeric
#    stoichiometric coefficients -1 / 1 that they should contain
#    (consistent with the "loss" reaction row which already used
#    numbers). Replace them with the correct numeric values.

$wb = $excel.ActiveWorkbook

$parsSheet = $wb.Worksheets.Item("pars")
$parsSheet.Range("D3").Value = 0.2
$parsSheet.Range("D4").Value = 0.4

$eqnsSheet = $wb.Worksheets.Item("eqns")
# transmission: S -1, E +1
$eqnsSheet.Range("E2").Value = -1
$eqnsSheet.Range("F2").Value = 1
# infection: E -1, I +1
$eqnsSheet.Range("F3").Value = -1
$eqnsSheet.Range("G3").Value = 1
# recovery: I -1, R +1
$eqnsSheet.Range("G4").Value = -1
$eqnsSheet.Range("H4").Value = 1

# Update the selected cell on the "pars" sheet (the active sheet) to
# match the author's final cursor position.
$parsSheet.Activate()
$parsSheet.Range("F13").Select()
